$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Step 1: Insert a fresh 12-row block after row 204 (before the old row 205),
# seeded with the *original* (pre-edit) content/format of rows 193:204 -- this
# captures the "skin_temperature_2" / unheated-window block exactly as it
# used to read, before we repurpose rows 181 and 193 below.
# ---------------------------------------------------------------------------
$ws.Range("A205:A216").EntireRow.Insert()
$ws.Range("A205:C216").RowHeight = 12

# Copy formats from the untouched source block (193:203) onto the new block.
$ws.Range("A193").Copy()
$ws.Range("A205").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C193").Copy()
$ws.Range("C205").PasteSpecial(-4122)

$ws.Range("C194:C203").Copy()
$ws.Range("C206:C215").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Values for the new block (identical text/format to the original 193:203
# "skin_temperature_2" definition).
$ws.Range("A205").Value = "skin_temperature_2"

$ws.Range("B206").Value = "type"
$ws.Range("C206").Value = "float32"

$ws.Range("B207").Value = "dimension"
$ws.Range("C207").Value = "time"

$ws.Range("B208").Value = "units"
$ws.Range("C208").Value = "K"

$ws.Range("B209").Value = "long_name"
$ws.Range("C209").Value = "Surface brightness temperature (9.6-11.5 um) at the Radiometer Stand (unheated window)"

$ws.Range("B210").Value = "standard_name"
$ws.Range("C210").Value = "surface_brightness_temperature"

$ws.Range("B211").Value = "_FillValue"
$ws.Range("C211").Value = -1.0 * [math]::Pow(10, 20)

$ws.Range("B212").Value = "cell_methods"
$ws.Range("C212").Value = "time: mean"

$ws.Range("B213").Value = "coordinates"
$ws.Range("C213").Value = "latitude longitude"

$ws.Range("B214").Value = "valid_min"
$ws.Range("C214").Value = "<derived from file>"

$ws.Range("B215").Value = "valid_max"
$ws.Range("C215").Value = "<derived from file>"

# ---------------------------------------------------------------------------
# Step 2: Relabel the two existing metadata blocks (rows 181:191 and
# 193:203) to reflect the new KT15 bias-correction variables.
# ---------------------------------------------------------------------------

# Met Mast (heated window) sensor -> now the *raw*, uncorrected variable.
$ws.Range("A181").Value = "skin_temperature_1_raw"
$ws.Range("C185").Value = "Surface brightness temperature (9.6-11.5 um) at the Met Mast (heated window), no bias correction applied"

# Former skin_temperature_2 block -> becomes skin_temperature_1, the
# bias-corrected Met Mast (heated window) variable.
$ws.Range("A193").Value = "skin_temperature_1"
$ws.Range("C197").Value = "Surface brightness temperature (9.6-11.5 um) at the Met Mast (heated window) with bias correction applied"

# ---------------------------------------------------------------------------
# Step 3: Update the sheet view (scroll position / selection) to match.
# ---------------------------------------------------------------------------
$ws.Range("C198").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 161
